# Updated cryptos list (Price + Volume(1h)) per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to stay text even when the new value looks numeric
    # (e.g. "577.36"), matching the existing inline-string cell type,
    # then drop the temporary text format so the cell style is untouched.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "64.850.10"
$ws.Range("E2").Value = "  -0.12%  "
Set-TextValue $ws.Range("D3") "3.155.69"
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("E4").Value = "  -0.12%  "
Set-TextValue $ws.Range("D5") "577.36"
$ws.Range("E5").Value = "  +1.21%  "
Set-TextValue $ws.Range("D6") "149.11"
$ws.Range("E6").Value = "  -0.87%  "
$ws.Range("E7").Value = "  -0.01%  "
Set-TextValue $ws.Range("D8") "3.153.41"
$ws.Range("E8").Value = "  +0.63%  "
Set-TextValue $ws.Range("D9") "0.525"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("E10").Value = "  -1.80%  "
Set-TextValue $ws.Range("D11") "6.11"
$ws.Range("E11").Value = "  -0.80%  "
Set-TextValue $ws.Range("D12") "0.500"
$ws.Range("E12").Value = "  -0.45%  "
Set-TextValue $ws.Range("D13") "0.0000260"
$ws.Range("E13").Value = "  +3.28%  "
Set-TextValue $ws.Range("D14") "37.14"
$ws.Range("E14").Value = "  -0.42%  "
Set-TextValue $ws.Range("D15") "3.668.28"
$ws.Range("E15").Value = "  +0.44%  "
Set-TextValue $ws.Range("D16") "64.961.64"
$ws.Range("E16").Value = "  +0.01%  "
Set-TextValue $ws.Range("D17") "3.148.59"
$ws.Range("E17").Value = "  +0.20%  "
Set-TextValue $ws.Range("D18") "7.11"
$ws.Range("E18").Value = "  -1.08%  "
$ws.Range("E19").Value = "  +0.39%  "
Set-TextValue $ws.Range("D20") "503.14"
$ws.Range("E20").Value = "  -1.16%  "
Set-TextValue $ws.Range("D21") "14.81"
$ws.Range("E21").Value = "  -0.72%  "
Set-TextValue $ws.Range("D22") "0.713"
$ws.Range("E22").Value = "  -2.46%  "
Set-TextValue $ws.Range("D23") "15.24"
$ws.Range("E23").Value = "  -1.75%  "
Set-TextValue $ws.Range("D24") "7.72"
$ws.Range("E24").Value = "  -1.20%  "
Set-TextValue $ws.Range("D25") "83.83"
$ws.Range("E25").Value = "  -1.55%  "
Set-TextValue $ws.Range("D26") "0.997"
$ws.Range("E26").Value = "  -0.20%  "
Set-TextValue $ws.Range("D27") "2.89"
$ws.Range("E27").Value = "  -0.89%  "
Set-TextValue $ws.Range("D28") "8.88"
$ws.Range("E28").Value = "  +2.10%  "
$ws.Range("E29").Value = "  -0.10%  "
Set-TextValue $ws.Range("D30") "2.83"
$ws.Range("E30").Value = "  +7.02%  "
Set-TextValue $ws.Range("D31") "27.54"
$ws.Range("E31").Value = "  -1.09%  "
Set-TextValue $ws.Range("D32") "0.999"
$ws.Range("E32").Value = "  -0.09%  "
Set-TextValue $ws.Range("D33") "1.20"
$ws.Range("E33").Value = "  +0.92%  "
Set-TextValue $ws.Range("D34") "6.17"
$ws.Range("E34").Value = "  +2.56%  "
Set-TextValue $ws.Range("D35") "6.46"
$ws.Range("E35").Value = "  -1.66%  "
Set-TextValue $ws.Range("D36") "54.56"
$ws.Range("E36").Value = "  -1.93%  "
Set-TextValue $ws.Range("D37") "0.0894"
$ws.Range("E37").Value = "  +4.26%  "
Set-TextValue $ws.Range("D38") "475.52"
$ws.Range("E38").Value = "  +0.61%  "
Set-TextValue $ws.Range("D39") "0.0415"
$ws.Range("E39").Value = "  -1.74%  "
Set-TextValue $ws.Range("D40") "2.97"
$ws.Range("E40").Value = "  -1.59%  "
Set-TextValue $ws.Range("D41") "8.65"
$ws.Range("E41").Value = "  +0.54%  "
Set-TextValue $ws.Range("D42") "3.010.84"
$ws.Range("E42").Value = "  -2.92%  "
$ws.Range("E43").Value = "  -3.87%  "
Set-TextValue $ws.Range("D44") "0.281"
$ws.Range("E44").Value = "  -2.93%  "
Set-TextValue $ws.Range("D45") "2.42"
$ws.Range("E45").Value = "  +0.00%  "
Set-TextValue $ws.Range("D46") "28.17"
$ws.Range("E46").Value = "  -3.18%  "
Set-TextValue $ws.Range("D47") "0.0₃0585"
$ws.Range("E47").Value = "  +1.80%  "
$ws.Range("E49").Value = "  -1.44%  "
Set-TextValue $ws.Range("D50") "2.23"
$ws.Range("E50").Value = "  -2.27%  "
Set-TextValue $ws.Range("D51") "33.57"
$ws.Range("E51").Value = "  +7.55%  "
